# Auto-generated: scheduled market-data refresh for Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18,8).Value = 653.6667
$ws.Cells.Item(18,9).Value = 653.6667
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,11).Value = 653.6667
$ws.Cells.Item(18,12).Value = 0
$ws.Cells.Item(18,13).Value = -369.6667
$ws.Cells.Item(18,14).ClearContents()
$ws.Cells.Item(33,8).Value = 6250105.5
$ws.Cells.Item(33,9).Value = 8333424
$ws.Cells.Item(33,10).Value = 150.5
$ws.Cells.Item(33,11).Value = 8333424
$ws.Cells.Item(33,12).Value = 150.5
$ws.Cells.Item(33,13).Value = -8333195
$ws.Cells.Item(33,14).Value = -608.5
$ws.Cells.Item(64,8).Value = 31252528
$ws.Cells.Item(64,9).Value = 50002012
$ws.Cells.Item(64,10).Value = 3387.5
$ws.Cells.Item(64,11).Value = 50002012
$ws.Cells.Item(64,12).Value = 3387.5
$ws.Cells.Item(64,13).Value = -50001764
$ws.Cells.Item(64,14).Value = -3883.5
$ws.Cells.Item(67,8).Value = 31252528
$ws.Cells.Item(67,9).Value = 50002012
$ws.Cells.Item(67,10).Value = 3387.5
$ws.Cells.Item(67,11).Value = 50002012
$ws.Cells.Item(67,12).Value = 3387.5
$ws.Cells.Item(67,13).Value = -50001154
$ws.Cells.Item(67,14).Value = -5103.5
$ws.Cells.Item(74,8).Value = 3218.423
$ws.Cells.Item(74,9).Value = 2517.9
$ws.Cells.Item(74,10).Value = 3656.25
$ws.Cells.Item(74,11).Value = 2517.9
$ws.Cells.Item(74,12).Value = 3656.25
$ws.Cells.Item(74,13).Value = -1581.9
$ws.Cells.Item(74,14).Value = -5528.25
$ws.Cells.Item(77,8).Value = 3218.423
$ws.Cells.Item(77,9).Value = 2517.9
$ws.Cells.Item(77,10).Value = 3656.25
$ws.Cells.Item(77,11).Value = 12589.5
$ws.Cells.Item(77,12).Value = 18281.25
$ws.Cells.Item(77,13).Value = -7909.5
$ws.Cells.Item(77,14).Value = -27641.25
$ws.Cells.Item(86,8).Value = 4329.2144
$ws.Cells.Item(86,9).Value = 2581.4443
$ws.Cells.Item(86,10).Value = 5157.1055
$ws.Cells.Item(86,11).Value = 2581.4443
$ws.Cells.Item(86,12).Value = 5157.1055
$ws.Cells.Item(86,13).Value = -1458.4443
$ws.Cells.Item(86,14).Value = -7403.1055
$ws.Cells.Item(88,8).Value = 9216.083000000001
$ws.Cells.Item(88,9).Value = 501.33334
$ws.Cells.Item(88,10).Value = 12121
$ws.Cells.Item(88,11).Value = 501.33334
$ws.Cells.Item(88,12).Value = 12121
$ws.Cells.Item(88,13).Value = -95.33334000000002
$ws.Cells.Item(88,14).Value = -12933
$ws.Cells.Item(89,8).Value = 4329.2144
$ws.Cells.Item(89,9).Value = 2581.4443
$ws.Cells.Item(89,10).Value = 5157.1055
$ws.Cells.Item(89,11).Value = 12907.2215
$ws.Cells.Item(89,12).Value = 25785.5275
$ws.Cells.Item(89,13).Value = -7291.2215
$ws.Cells.Item(89,14).Value = -37017.5275
$ws.Cells.Item(91,8).Value = 9216.083000000001
$ws.Cells.Item(91,9).Value = 501.33334
$ws.Cells.Item(91,10).Value = 12121
$ws.Cells.Item(91,11).Value = 501.33334
$ws.Cells.Item(91,12).Value = 12121
$ws.Cells.Item(91,13).Value = 902.66666
$ws.Cells.Item(91,14).Value = -14929
$ws.Cells.Item(118,8).Value = 885.09753
$ws.Cells.Item(118,9).Value = 533.3333
$ws.Cells.Item(118,11).Value = 1599.9999
$ws.Cells.Item(118,13).Value = 57.00009999999997
$ws.Cells.Item(137,8).Value = 4864.9697
$ws.Cells.Item(137,9).Value = 3235
$ws.Cells.Item(137,10).Value = 5227.185
$ws.Cells.Item(137,11).Value = 9705
$ws.Cells.Item(137,12).Value = 15681.555
$ws.Cells.Item(137,13).Value = -7155
$ws.Cells.Item(137,14).Value = -20781.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24,8).Value = 26451.666
$ws.Cells.Item(24,10).Value = 26451.666
$ws.Cells.Item(24,12).Value = 26451.666
$ws.Cells.Item(24,14).Value = -27199.666
$ws.Cells.Item(28,8).Value = 32333.334
$ws.Cells.Item(28,9).Value = 9000
$ws.Cells.Item(28,10).Value = 44000
$ws.Cells.Item(28,11).Value = 9000
$ws.Cells.Item(28,12).Value = 44000
$ws.Cells.Item(28,13).Value = -8808
$ws.Cells.Item(28,14).Value = -44384
$ws.Cells.Item(32,8).Value = 7059.415
$ws.Cells.Item(32,9).Value = 6571.549
$ws.Cells.Item(32,11).Value = 6571.549
$ws.Cells.Item(32,13).Value = -6284.549
$ws.Cells.Item(93,8).Value = 18950
$ws.Cells.Item(93,10).Value = 18950
$ws.Cells.Item(93,12).Value = 18950
$ws.Cells.Item(93,14).Value = -23942
$ws.Cells.Item(94,8).Value = 22500
$ws.Cells.Item(94,10).Value = 22500
$ws.Cells.Item(94,12).Value = 22500
$ws.Cells.Item(94,14).Value = -24302
$ws.Cells.Item(95,8).Value = 15801.6
$ws.Cells.Item(95,10).Value = 15801.6
$ws.Cells.Item(95,12).Value = 15801.6
$ws.Cells.Item(95,14).Value = -21293.6
$ws.Cells.Item(96,8).Value = 40000
$ws.Cells.Item(96,10).Value = 40000
$ws.Cells.Item(96,12).Value = 40000
$ws.Cells.Item(96,14).Value = -45492
$ws.Cells.Item(98,8).Value = 48000
$ws.Cells.Item(98,10).Value = 48000
$ws.Cells.Item(98,12).Value = 48000
$ws.Cells.Item(98,14).Value = -53990
$ws.Cells.Item(99,8).Value = 32333.334
$ws.Cells.Item(99,9).Value = 9000
$ws.Cells.Item(99,10).Value = 44000
$ws.Cells.Item(99,11).Value = 9000
$ws.Cells.Item(99,12).Value = 44000
$ws.Cells.Item(99,13).Value = -6005
$ws.Cells.Item(99,14).Value = -49990
$ws.Cells.Item(100,8).Value = 26451.666
$ws.Cells.Item(100,10).Value = 26451.666
$ws.Cells.Item(100,12).Value = 26451.666
$ws.Cells.Item(100,14).Value = -28615.666
$ws.Cells.Item(101,8).Value = 48000
$ws.Cells.Item(101,10).Value = 48000
$ws.Cells.Item(101,12).Value = 48000
$ws.Cells.Item(101,14).Value = -54490
$ws.Cells.Item(102,8).Value = 1142.1
$ws.Cells.Item(102,9).Value = 1038.75
$ws.Cells.Item(102,10).Value = 1555.5
$ws.Cells.Item(102,11).Value = 1038.75
$ws.Cells.Item(102,12).Value = 1555.5
$ws.Cells.Item(102,13).Value = 583.25
$ws.Cells.Item(102,14).Value = -4799.5
$ws.Cells.Item(103,8).Value = 0
$ws.Cells.Item(103,10).Value = 0
$ws.Cells.Item(103,12).Value = 0
$ws.Cells.Item(103,14).ClearContents()
$ws.Cells.Item(104,8).Value = 22747.5
$ws.Cells.Item(104,10).Value = 22747.5
$ws.Cells.Item(104,12).Value = 22747.5
$ws.Cells.Item(104,14).Value = -29735.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(104,8).Value = 32500
$ws.Cells.Item(104,10).Value = 32500
$ws.Cells.Item(104,12).Value = 32500
$ws.Cells.Item(104,14).Value = -39488
$ws.Cells.Item(105,8).Value = 2044.8
$ws.Cells.Item(105,9).Value = 1388.3334
$ws.Cells.Item(105,10).Value = 2482.4443
$ws.Cells.Item(105,11).Value = 1388.3334
$ws.Cells.Item(105,12).Value = 2482.4443
$ws.Cells.Item(105,13).Value = 358.6666
$ws.Cells.Item(105,14).Value = -5976.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 28301.791
$ws.Cells.Item(31,9).Value = 1620.6666
$ws.Cells.Item(31,10).Value = 32628.459
$ws.Cells.Item(31,11).Value = 1620.6666
$ws.Cells.Item(31,12).Value = 32628.459
$ws.Cells.Item(31,13).Value = -1325.6666
$ws.Cells.Item(31,14).Value = -33218.459
$ws.Cells.Item(34,8).Value = 28301.791
$ws.Cells.Item(34,9).Value = 1620.6666
$ws.Cells.Item(34,10).Value = 32628.459
$ws.Cells.Item(34,11).Value = 1620.6666
$ws.Cells.Item(34,12).Value = 32628.459
$ws.Cells.Item(34,13).Value = -1418.6666
$ws.Cells.Item(34,14).Value = -33032.459
$ws.Cells.Item(99,8).Value = 2249.1904
$ws.Cells.Item(99,9).Value = 1451.5
$ws.Cells.Item(99,10).Value = 2974.3635
$ws.Cells.Item(99,11).Value = 1451.5
$ws.Cells.Item(99,12).Value = 2974.3635
$ws.Cells.Item(99,13).Value = 46.5
$ws.Cells.Item(99,14).Value = -5970.363499999999
$ws.Cells.Item(126,8).Value = 2249.1904
$ws.Cells.Item(126,9).Value = 1451.5
$ws.Cells.Item(126,10).Value = 2974.3635
$ws.Cells.Item(126,11).Value = 4354.5
$ws.Cells.Item(126,12).Value = 8923.0905
$ws.Cells.Item(126,13).Value = -1884.5
$ws.Cells.Item(126,14).Value = -13863.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131,8).Value = 929
$ws.Cells.Item(131,9).Value = 929
$ws.Cells.Item(131,10).Value = 0
$ws.Cells.Item(131,11).Value = 2787
$ws.Cells.Item(131,12).Value = 0
$ws.Cells.Item(131,13).Value = 2253
$ws.Cells.Item(131,14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122,8).Value = 843.1429000000001
$ws.Cells.Item(122,9).Value = 878.63635
$ws.Cells.Item(122,10).Value = 713
$ws.Cells.Item(122,11).Value = 2635.90905
$ws.Cells.Item(122,12).Value = 2139
$ws.Cells.Item(122,13).Value = -185.9090500000002
$ws.Cells.Item(122,14).Value = -7039

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62,8).Value = 2682.077
$ws.Cells.Item(62,9).Value = 2727.8333
$ws.Cells.Item(62,10).Value = 2642.8572
$ws.Cells.Item(62,11).Value = 2727.8333
$ws.Cells.Item(62,12).Value = 2642.8572
$ws.Cells.Item(62,13).Value = -2103.8333
$ws.Cells.Item(62,14).Value = -3890.8572
$ws.Cells.Item(65,8).Value = 2682.077
$ws.Cells.Item(65,9).Value = 2727.8333
$ws.Cells.Item(65,10).Value = 2642.8572
$ws.Cells.Item(65,11).Value = 13639.1665
$ws.Cells.Item(65,12).Value = 13214.286
$ws.Cells.Item(65,13).Value = -10519.1665
$ws.Cells.Item(65,14).Value = -19454.286
